# Update gh-pages to output generated at 456a3b4
# Increment the "想去人数" (want-to-go count) in column F for two events
# that appear on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5403
$ws1.Range("F12").Value = 36

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 5403
$ws4.Range("F14").Value = 36
